$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.610.00'
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").Value = '3.509.97'
$ws.Range("E3").Value = '  -0.42%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.27'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.59%  '
$ws.Range("D7").Value = '3.509.28'
$ws.Range("E7").Value = '  -0.41%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.488'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.07%  '
$ws.Range("E10").Value = '  +2.83%  '
$ws.Range("E11").Value = '  +7.63%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.434'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.95%  '
$ws.Range("E13").Value = '  -1.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.25'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.54%  '
$ws.Range("D15").Value = '4.101.30'
$ws.Range("E15").Value = '  -0.46%  '
$ws.Range("D16").Value = '3.511.71'
$ws.Range("E16").Value = '  -0.49%  '
$ws.Range("D17").Value = '67.501.17'
$ws.Range("E17").Value = '  +0.12%  '
$ws.Range("E18").Value = '  -0.65%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.52'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.09%  '
$ws.Range("E20").Value = '  +1.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.85'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.27%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '447.30'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.55%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.630'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.87%  '
$ws.Range("E24").Value = '  +1.07%  '
$ws.Range("D25").Value = '3.649.07'
$ws.Range("E25").Value = '  -0.45%  '
$ws.Range("E26").Value = '  -2.90%  '
$ws.Range("E27").Value = '  +0.03%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.75'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.40%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.06'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.79%  '
$ws.Range("E30").Value = '  +0.36%  '
$ws.Range("E31").Value = '  +5.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.174'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.02%  '
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("E34").Value = '  -0.88%  '
$ws.Range("E35").Value = '  -0.27%  '
$ws.Range("E36").Value = '  +1.06%  '
$ws.Range("D37").Value = '3.498.36'
$ws.Range("E37").Value = '  -0.49%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.99'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.12%  '
$ws.Range("E40").Value = '  +7.22%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '179.08'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0897'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.57%  '
$ws.Range("E44").Value = '  +0.54%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.893'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.45%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '30.13'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.96%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '46.46'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.93%  '
$ws.Range("E48").Value = '  +4.39%  '
$ws.Range("E49").Value = '  -3.19%  '
$ws.Range("E50").Value = '  +0.39%  '
$ws.Range("E51").Value = '  +1.97%  '
